{"js": "// Update the Chai Tea \"lemas\" (slogans) list so each bullet carries the\n// \"T\u00e9 Chai:\" brand prefix, per the commit's canonical diff.\nconst replacements = [\n  [\"un mundo de sabor en una taza\", \"T\u00e9 Chai: un mundo de sabor en una taza\"],\n  [\"la combinaci\u00f3n perfecta de salud y placer\", \"T\u00e9 Chai: La mezcla perfecta de salud y placer\"],\n  [\"m\u00e1s que un simple t\u00e9, una forma de vida\", \"T\u00e9 Chai: M\u00e1s que solo t\u00e9, una forma de vida\"],\n  [\"una bebida para cualquier estaci\u00f3n y motivo\", \"T\u00e9 Chai: Una bebida para todas las estaciones y razones\"],\n  [\"el placer definitivo para sus sentidos\", \"T\u00e9 Chai: la m\u00e1xima indulgencia para sus sentidos\"],\n  [\"una dulce v\u00eda de escape de la rutina\", \"T\u00e9 Chai: Un dulce escape del d\u00eda a d\u00eda\"],\n  [\"comparta el calor y el amor\", \"T\u00e9 Chai: Compartir la calidez, compartir el amor\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n\n  results.items[0].insertText(newText, \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "# Update the Chai Tea \"lemas\" (slogans) list so each bullet carries the\n# \"T\u00e9 Chai:\" brand prefix, per the commit's canonical diff.\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"un mundo de sabor en una taza\", \"T\u00e9 Chai: un mundo de sabor en una taza\"),\n  @(\"la combinaci\u00f3n perfecta de salud y placer\", \"T\u00e9 Chai: La mezcla perfecta de salud y placer\"),\n  @(\"m\u00e1s que un simple t\u00e9, una forma de vida\", \"T\u00e9 Chai: M\u00e1s que solo t\u00e9, una forma de vida\"),\n  @(\"una bebida para cualquier estaci\u00f3n y motivo\", \"T\u00e9 Chai: Una bebida para todas las estaciones y razones\"),\n  @(\"el placer definitivo para sus sentidos\", \"T\u00e9 Chai: la m\u00e1xima indulgencia para sus sentidos\"),\n  @(\"una dulce v\u00eda de escape de la rutina\", \"T\u00e9 Chai: Un dulce escape del d\u00eda a d\u00eda\"),\n  @(\"comparta el calor y el amor\", \"T\u00e9 Chai: Compartir la calidez, compartir el amor\")\n)\n\nforeach ($pair in $pairs) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $found = $find.Execute($oldText, $true, $true, $false, $false, $false, $true, 1, $false, $newText, 2)\n\n  if (-not $found) {\n    Write-Output \"NOT FOUND: $oldText\"\n  }\n}\n"}
